# Resolves search issue by restricting search.jsp. Stops user from logging
# in while attempting search, which was what messed things up before.
#
# Inserts 3 new rows above row 55 (pushing the existing "Thursday plan…"
# block and everything below it down by 3 rows) and fills the two new
# rows with a note about web.xml plus a reference link, leaving the third
# new row blank (matches the existing blank-row pattern used elsewhere in
# the log).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows immediately above row 55.
$ws.Rows("55:57").Insert() | Out-Null

# New note + reference link in the freshly inserted rows.
$ws.Range("D55").Value = "Note:  I may have needed to do something in web xml to get my application startup servlet to actually load on startup"
$ws.Range("D56").Value = "http://tutorials.jenkov.com/java-servlets/web-xml.html"

# Update the view to match the author's saved state (scroll + final selection).
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D52").Select() | Out-Null
